$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="29.475.11"'
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("D3").Formula = '="1.840.65"'
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("D4").Formula = '="0.9983"'
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("E4").Value = '  -1.07%  '
$ws.Range("D5").Formula = '="244.34"'
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").Formula = '="0.6278"'
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("D7").Formula = '="0.9992"'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("D8").Formula = '="0.07439"'
$ws.Range("D8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("D9").Formula = '="0.2950"'
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = '  +1.50%  '
$ws.Range("D10").Formula = '="23.65"'
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Value = '  +3.39%  '
$ws.Range("D11").Formula = '="0.07654"'
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").Formula = '="1.834.02"'
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Formula = '="5.019"'
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("D14").Formula = '="0.6778"'
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = '  +1.55%  '
$ws.Range("D15").Formula = '="83.77"'
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Formula = '="0.000009360"'
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$ws.Range("E16").Value = '  +3.27%  '
$ws.Range("D17").Formula = '="5.929"'
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Value = '  +1.08%  '
$ws.Range("D18").Formula = '="29.445.65"'
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").Formula = '="2.085.67"'
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Formula = '="238.05"'
$ws.Range("D20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("D23").Formula = '="7.358"'
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Value = '  +3.05%  '
$ws.Range("D24").Formula = '="0.9998"'
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("D25").Formula = '="159.14"'
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Formula = '="0.1414"'
$ws.Range("D26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Formula = '="8.521"'
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = '  +0.53%  '
$ws.Range("D28").Formula = '="17.76"'
$ws.Range("D28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").Formula = '="0.06037"'
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("E30").Value = '  +8.64%  '
$ws.Range("D31").Formula = '="1.244"'
$ws.Range("D31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = '  +2.37%  '
$ws.Range("D32").Formula = '="4.101"'
$ws.Range("D32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").Formula = '="4.117"'
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").Formula = '="1.876"'
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("D35").Formula = '="1.143"'
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = '  +0.71%  '
$ws.Range("D36").Formula = '="0.7281"'
$ws.Range("D36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = '  -1.53%  '
$ws.Range("D37").Formula = '="2.613"'
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = '  -1.49%  '
$ws.Range("D38").Formula = '="2.879"'
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = '  +2.04%  '
$ws.Range("D39").Formula = '="1.220.80"'
$ws.Range("D39").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$ws.Range("E39").Value = '  +1.53%  '
$ws.Range("D40").Formula = '="0.01764"'
$ws.Range("D40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = '  -0.45%  '
$ws.Range("D41").Formula = '="6.272"'
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("E41").Value = '  -2.07%  '
$ws.Range("D42").Formula = '="0.9127"'
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("D43").Formula = '="1.000"'
$ws.Range("D43").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("D44").Formula = '="2.001.90"'
$ws.Range("D44").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("D45").Formula = '="102.02"'
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("D46").Formula = '="65.58"'
$ws.Range("D46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Value = '  +1.33%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Formula = '="0.5069"'
$ws.Range("D47").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Formula = '="9.275"'
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = '  +1.76%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Formula = '="0.00000000118"'
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").Formula = '="0.4061"'
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").Formula = '="0.1139"'
$ws.Range("D51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = '  +3.49%  '
$excel.CutCopyMode = 0
